# Update "Resumo Inscricoes" counts (Inscritos/Pagos/Inscricoes homologadas)
# on the Inscricoes sheet, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F2").Value = 37
$ws.Range("H2").Value = 37

$ws.Range("E4").Value = 27

$ws.Range("E5").Value = 66

$ws.Range("F10").Value = 86
$ws.Range("H10").Value = 86

$ws.Range("E11").Value = 145
$ws.Range("F11").Value = 76
$ws.Range("H11").Value = 76

$ws.Range("E12").Value = 219

$ws.Range("E17").Value = 41

$ws.Range("E21").Value = 70
$ws.Range("F21").Value = 36
$ws.Range("H21").Value = 36

$ws.Range("E23").Value = 97
$ws.Range("F23").Value = 41
$ws.Range("H23").Value = 41

$ws.Range("E25").Value = 94

$ws.Range("E26").Value = 56

$ws.Range("E27").Value = 141
$ws.Range("F27").Value = 68
$ws.Range("H27").Value = 68

$ws.Range("E28").Value = 90

$ws.Range("E29").Value = 91

$ws.Range("E30").Value = 100

$ws.Range("E33").Value = 129

$ws.Range("E34").Value = 98
$ws.Range("F34").Value = 53
$ws.Range("H34").Value = 53

$ws.Range("E35").Value = 64

$ws.Range("E37").Value = 67

$ws.Range("E40").Value = 139

$ws.Range("E41").Value = 176

$ws.Range("E44").Value = 138
$ws.Range("F44").Value = 65
$ws.Range("H44").Value = 65

$ws.Range("E46").Value = 121
$ws.Range("F46").Value = 56
$ws.Range("H46").Value = 56

$ws.Range("E47").Value = 200

$ws.Range("E51").Value = 97
